$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that currently sits at the end of the
#    third paragraph ("...tit for tat opponent.").
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Expand the sentence about identifying teammates / defecting against
#    non-teammates in the second paragraph, right before the existing
#    "The team player strategy sacrifices..." sentence.
$oldText = "everyone else as little points as possible.  The team player strategy sacrifices"
$newText = "everyone else as little points as possible.  " + `
    "First it identifies " + `
    "the opposing player as its teammate.  Once the teammate is identified, the team player strategy continuously cooperates while the teammate continuously defects.  If team player" + [char]0x2019 + "s opponent is identified as not a teammate, it will continuously defect " + `
    "in order to" + `
    " lower the opponent" + [char]0x2019 + "s overall score. " + `
    "The team player strategy sacrifices"

[void]$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# 3. Re-insert the "_GoBack" bookmark right after the newly added text, just
#    before "The team player strategy sacrifices...".
$marker = $d.Content
[void]$marker.Find.Execute("The team player strategy sacrifices itself to elevate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($marker.Start, $marker.Start)
[void]$d.Bookmarks.Add("_GoBack", $bmRange)
